$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.38222733333333
$ws.Range("H2").Value = 76.146682
$ws.Range("I2").Value = 0.1760862452187379
$ws.Range("J2").Value = 0.1760862452187379
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 2855.840493788448
$ws.Range("R2").Value = 25702.56444409603
$ws.Range("S2").Value = 0.05767375641010355
$ws.Range("T2").Value = 0.05767375641010355
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.38222733333333
$ws.Range("H3").Value = 76.146682
$ws.Range("I3").Value = 0.1760862452187379
$ws.Range("J3").Value = 0.1760862452187379
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 2698.497944833937
$ws.Range("R3").Value = 24286.48150350543
$ws.Range("S3").Value = 0.05449622045839873
$ws.Range("T3").Value = 0.05449622045839874
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.38222733333333
$ws.Range("H4").Value = 76.146682
$ws.Range("I4").Value = 0.1760862452187379
$ws.Range("J4").Value = 0.1760862452187379
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 3164.951942240313
$ws.Range("R4").Value = 28484.56748016281
$ws.Range("S4").Value = 0.06391626835023563
$ws.Range("T4").Value = 0.06391626835023566
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 94.773687
$ws.Range("H5").Value = 284.321061
$ws.Range("I5").Value = 0.6574814128880592
$ws.Range("J5").Value = 0.6574814128880593
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 10663.3089967163
$ws.Range("R5").Value = 95969.78097044672
$ws.Range("S5").Value = 0.2153457404010879
$ws.Range("T5").Value = 0.2153457404010879
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 94.773687
$ws.Range("H6").Value = 284.321061
$ws.Range("I6").Value = 0.6574814128880592
$ws.Range("J6").Value = 0.6574814128880593
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 10075.81392425614
$ws.Range("R6").Value = 90682.32531830527
$ws.Range("S6").Value = 0.2034812655556263
$ws.Range("T6").Value = 0.2034812655556264
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 94.773687
$ws.Range("H7").Value = 284.321061
$ws.Range("I7").Value = 0.6574814128880592
$ws.Range("J7").Value = 0.6574814128880593
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 11817.48791407321
$ws.Range("R7").Value = 106357.3912266589
$ws.Range("S7").Value = 0.2386544069313449
$ws.Range("T7").Value = 0.238654406931345
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.99065033333333
$ws.Range("H8").Value = 71.971951
$ws.Range("I8").Value = 0.1664323418932028
$ws.Range("J8").Value = 0.1664323418932028
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 2699.269445289264
$ws.Range("R8").Value = 24293.42500760337
$ws.Range("S8").Value = 0.05451180092566488
$ws.Range("T8").Value = 0.05451180092566488
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.99065033333333
$ws.Range("H9").Value = 71.971951
$ws.Range("I9").Value = 0.1664323418932028
$ws.Range("J9").Value = 0.1664323418932028
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 2550.553179181055
$ws.Range("R9").Value = 22954.9786126295
$ws.Range("S9").Value = 0.05150847293014121
$ws.Range("T9").Value = 0.05150847293014123
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.99065033333333
$ws.Range("H10").Value = 71.971951
$ws.Range("I10").Value = 0.1664323418932028
$ws.Range("J10").Value = 0.1664323418932028
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 2991.433902586519
$ws.Range("R10").Value = 26922.90512327867
$ws.Range("S10").Value = 0.0604120680373967
$ws.Range("T10").Value = 0.06041206803739672
Write-Output "applied changes"
